$wb = $excel.ActiveWorkbook
$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# --- Append three new rows to the "Logs" sheet ---

$wsLogs.Cells.Item(33, 1).Value = "Vragen over samenwerking"
$wsLogs.Cells.Item(33, 2).Value = "mailmind.test@zohomail.eu"
$wsLogs.Cells.Item(33, 3).Value = "Kunnen we samenwerken aan een nieuw project?"
$wsLogs.Cells.Item(33, 4).Value = "Samenwerking / Partnerverzoek"
$wsLogs.Cells.Item(33, 6).Value = "2025-06-19 22:06:32"
$wsLogs.Cells.Item(33, 7).Value = "Nee"

$wsLogs.Cells.Item(34, 1).Value = "Klacht over levering"
$wsLogs.Cells.Item(34, 2).Value = "mailmind.test@zohomail.eu"
$wsLogs.Cells.Item(34, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$wsLogs.Cells.Item(34, 4).Value = "Klacht / Probleem"
$wsLogs.Cells.Item(34, 6).Value = "2025-06-19 22:06:33"
$wsLogs.Cells.Item(34, 7).Value = "Nee"

$wsLogs.Cells.Item(35, 1).Value = "Uitnodiging voor netwerkevent"
$wsLogs.Cells.Item(35, 2).Value = "mailmind.test@zohomail.eu"
$wsLogs.Cells.Item(35, 3).Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$wsLogs.Cells.Item(35, 4).Value = "Samenwerking / Partnerverzoek"
$wsLogs.Cells.Item(35, 6).Value = "2025-06-19 22:06:33"
$wsLogs.Cells.Item(35, 7).Value = "Nee"

# --- Extend the conditional formatting ranges to cover the new rows ---

$fcCategory = $wsLogs.Range("D2:D32").FormatConditions.Item(1)
$fcCategory.ModifyAppliesToRange($wsLogs.Range("D2:D35"))

$fcAnswered = $wsLogs.Range("G2:G32").FormatConditions.Item(1)
$fcAnswered.ModifyAppliesToRange($wsLogs.Range("G2:G35"))

# --- Update the "Dashboard" summary sheet ---

# Samenwerking / Partnerverzoek count: 8 -> 10
$wsDash.Cells.Item(2, 2).Value = 10

# Rows 7 and 8 swap category order, with updated counts
$wsDash.Cells.Item(7, 1).Value = "Klacht / Probleem"
$wsDash.Cells.Item(7, 2).Value = 3
$wsDash.Cells.Item(8, 1).Value = "Offerte / Prijsaanvraag"
$wsDash.Cells.Item(8, 2).Value = 3
